# Horarios actualizados Linea 141 - 678
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with a
# fresh scrape timestamp (03:22:47 -> 03:58:57) and new/updated rows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # LP1912
$ws2 = $wb.Worksheets.Item(2)   # LP1912-215
$ws3 = $wb.Worksheets.Item(3)   # 6203-6173

$newTime = "03:58:57"

# ---------------------------------------------------------------------
# Sheet 1: LP1912 - a bus already arrived (14_ABASTO dropped), every row
# shifts up one position, refreshed with three new rows at the bottom.
# ---------------------------------------------------------------------
$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 8"

$ws1.Cells.Item(6,1).Value = $newTime
$ws1.Cells.Item(6,2).Value = "04:01"
$ws1.Cells.Item(6,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(6,4).Value = 3
$ws1.Cells.Item(6,5).Value = "LP1912"

$ws1.Cells.Item(7,1).Value = $newTime
$ws1.Cells.Item(7,2).Value = "04:45"
$ws1.Cells.Item(7,3).Value = "215A_EL PATO"
$ws1.Cells.Item(7,4).Value = 47
$ws1.Cells.Item(7,5).Value = "LP1912"

$ws1.Cells.Item(8,1).Value = $newTime
$ws1.Cells.Item(8,2).Value = "04:53"
$ws1.Cells.Item(8,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(8,4).Value = 55
$ws1.Cells.Item(8,5).Value = "LP1912"

$ws1.Cells.Item(9,1).Value = $newTime
$ws1.Cells.Item(9,2).Value = "05:16"
$ws1.Cells.Item(9,3).Value = "17_ROMERO"
$ws1.Cells.Item(9,4).Value = 78
$ws1.Cells.Item(9,5).Value = "LP1912"

$ws1.Cells.Item(10,1).Value = $newTime
$ws1.Cells.Item(10,2).Value = "05:21"
$ws1.Cells.Item(10,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(10,4).Value = 83
$ws1.Cells.Item(10,5).Value = "LP1912"

$ws1.Cells.Item(11,1).Value = $newTime
$ws1.Cells.Item(11,2).Value = "05:34"
$ws1.Cells.Item(11,3).Value = "215B_EL PATO"
$ws1.Cells.Item(11,4).Value = 96
$ws1.Cells.Item(11,5).Value = "LP1912"

$ws1.Cells.Item(12,1).Value = $newTime
$ws1.Cells.Item(12,2).Value = "05:46"
$ws1.Cells.Item(12,3).Value = "15_ABASTO"
$ws1.Cells.Item(12,4).Value = 108
$ws1.Cells.Item(12,5).Value = "LP1912"

$ws1.Cells.Item(13,1).Value = $newTime
$ws1.Cells.Item(13,2).Value = "05:53"
$ws1.Cells.Item(13,3).Value = "10_OLMOS"
$ws1.Cells.Item(13,4).Value = 115
$ws1.Cells.Item(13,5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215 - refreshed existing row, plus a new 215B arrival.
# ---------------------------------------------------------------------
$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2.Cells.Item(6,1).Value = $newTime
$ws2.Cells.Item(6,2).Value = "04:45"
$ws2.Cells.Item(6,3).Value = "215A_EL PATO"
$ws2.Cells.Item(6,4).Value = 47
$ws2.Cells.Item(6,5).Value = "LP1912"

$ws2.Cells.Item(7,1).Value = $newTime
$ws2.Cells.Item(7,2).Value = "05:34"
$ws2.Cells.Item(7,3).Value = "215B_EL PATO"
$ws2.Cells.Item(7,4).Value = 96
$ws2.Cells.Item(7,5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173 - previously empty (0 rows); now has its first
# tracked arrival, so the header row needs to be (re)created too. Copy
# it from sheet 1's header so the styling (bold, border, alignment)
# matches exactly.
# ---------------------------------------------------------------------
$ws3.Range("A2").Value = "Última actualización: $newTime"
$ws3.Range("A3").Value = "Total filas: 1"

$ws1.Range("A5:E5").Copy($ws3.Range("A5:E5"))

$ws3.Cells.Item(6,1).Value = $newTime
$ws3.Cells.Item(6,2).Value = "05:43"
$ws3.Cells.Item(6,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(6,4).Value = 105
$ws3.Cells.Item(6,5).Value = "L6173"
